# Remove needless imports on org.dozer
#
# The "Environment" / "import" table on Sheet1 lists imported packages in
# column D, rows 8-12 (merged C8:C12). The "org.dozer" import (row 12) is no
# longer needed, so the entire row is removed, shifting everything below it
# up by one row (merge C8:C12 -> C8:C11, data table moves from rows 15-20 to
# rows 14-19, etc.). Excel/ the workbook also drops the now-unused
# "org.dozer" shared string automatically.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Delete the whole row that holds the "org.dozer" import.
$ws1.Rows("12:12").Delete()

# Restore a sensible current selection on the sheet (author ended up with
# D12 selected, which now is the last row of the shrunk import list).
$ws1.Range("D12").Select() | Out-Null
